# Apply the edits described by the commit:
# "Trade #10 closed at 2026-02-17 12:27:36 - unknown UNKNOWN +0.000%"
#
# 1) Update aggregate stats on the "Summary" sheet
# 2) Update aggregate stats for the "MarketMaking" strategy row on the
#    "Strategy Status" sheet
# 3) Append the new (10th) closed trade as row 11 on both the
#    "All Trades" and "MarketMaking" trade-log sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.02   # Current Capital
$summary.Range("B4").Value = 0.02      # Total P&L $
$summary.Range("B5").Value = 0.04      # Total P&L %
$summary.Range("B6").Value = 10        # Total Trades
$summary.Range("B8").Value = 4         # Losing Trades
$summary.Range("B9").Value = 40        # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.02     # Capital
$status.Range("D4").Value = 10         # Trades
$status.Range("E4").Value = 0.02       # P&L $
$status.Range("F4").Value = 0.02       # P&L %
$status.Range("G4").Value = 40         # Win Rate %

# ---------------------------------------------------------------------
# 3) Append trade #10 to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 11

    $ws.Cells.Item($row, 1).Value = 10          # Trade #

    # Text-like values that would otherwise be auto-converted to
    # dates/times/numbers by Excel's smart entry need to be forced to
    # Text format first, then the style reset back to Normal so no
    # stray cell style is left behind.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17" # Date
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "12:27:29"   # Time
    $ws.Cells.Item($row, 3).Style = "Normal"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"  # Strategy
    $ws.Cells.Item($row, 5).Value = "UP"            # Side
    $ws.Cells.Item($row, 6).Value = 0.09            # Entry Price
    $ws.Cells.Item($row, 7).Value = 0.07000000000000001 # Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"        # Status
    $ws.Cells.Item($row, 9).Value = -22.2222        # P&L %
    $ws.Cells.Item($row, 10).Value = -0.02          # P&L $
    $ws.Cells.Item($row, 11).Value = 100.02         # Capital After
    $ws.Cells.Item($row, 12).Value = 0              # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0              # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6            # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps" # Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"   # Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.13           # Duration (min)
}
